$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------
# 1. Update the data table (rows 2-10): Date / Arrivals / Closed
#    Columns D:G are formulas and recalc automatically.
# ---------------------------------------------------------------
$ws.Range("A2").Value = "7/1/2011"
$ws.Range("B2").Value = 24
$ws.Range("C2").Value = 3

$ws.Range("A3").Value = "8/1/2011"
$ws.Range("B3").Value = 15
$ws.Range("C3").Value = 3

$ws.Range("A4").Value = "9/1/2011"
$ws.Range("B4").Value = 8
$ws.Range("C4").Value = 0

$ws.Range("A5").Value = "7/1/2017"
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = 0

$ws.Range("A6").Value = "8/1/2017"
$ws.Range("B6").Value = 5
$ws.Range("C6").Value = 4

$ws.Range("A7").Value = "9/1/2017"
$ws.Range("B7").Value = 4
$ws.Range("C7").Value = 2

$ws.Range("A8").Value = "1/1/2018"
$ws.Range("B8").Value = 5
$ws.Range("C8").Value = 2

$ws.Range("A9").Value = "2/1/2018"
$ws.Range("B9").Value = 3
$ws.Range("C9").Value = 1

$ws.Range("A10").Value = "3/1/2018"
$ws.Range("B10").Value = 7
$ws.Range("C10").Value = 2

# Keep the date column's number format (mmm-yy) intact for the rewritten cells.
$ws.Range("A2:A10").NumberFormat = "mmm-yy"

# ---------------------------------------------------------------
# 2. Version-coverage labels in column H (bucket legend per quarter)
# ---------------------------------------------------------------
$ws.Range("H2").Value = "3.0~3.6"
$ws.Range("H5").Value = "3.6~3.7"
$ws.Range("H8").Value = "3.7~3.8"

# ---------------------------------------------------------------
# 3. Update the JQL query text (row 12) and the report URL (row 13)
# ---------------------------------------------------------------
$ws.Range("A12").Value = "project = LANG AND issuetype = Bug AND  affectedVersion in (3.0,3.6) AND status = Closed  AND createdDate > ""2011/01/01"" ORDER BY created DESC"

$newUrl = "https://issues.apache.org/jira/secure/ConfigureReport.jspa?projectOrFilterId=project-12310481&dateField=created&periodName=monthly&daysprevious=7200&cumulative=true&selectedProjectId=12310481&reportKey=com.atlassian.jira.jira-core-reports-plugin%3Atimesince-report&atl_token=A5KQ-2QAV-T4JA-FDED%7C06b3f857d3c464eb55cc90414bee10b3a9f28b29%7Clout&Next=Next"

$ws.Range("A13").Value = $newUrl
$ws.Hyperlinks.Add($ws.Range("A13"), $newUrl, "", "", $newUrl) | Out-Null
$ws.Range("A13").Style = "超链接"

# Row 13 (besides A13) keeps plain/default formatting, same as row 12.
$ws.Range("B13:M13").Style = "常规"

# ---------------------------------------------------------------
# 4. Cosmetic: move the active selection (matches the saved view state)
# ---------------------------------------------------------------
$ws.Range("D17").Select() | Out-Null

Write-Output "done"
